$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.521.82"
$ws.Range("E2").Value = "  -1.99%  "

$ws.Range("D3").Value = "1.750.34"
$ws.Range("E3").Value = "  -2.22%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Formula = '="323.91"'
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E5").Value = "  +0.17%  "

$ws.Range("D6").Formula = '="0.9997"'
$ws.Range("D6").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").Formula = '="0.4467"'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E7").Value = "  +3.56%  "

$ws.Range("D8").Formula = '="0.3602"'
$ws.Range("D8").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E8").Value = "  -0.59%  "

$ws.Range("D9").Formula = '="0.07510"'
$ws.Range("D9").Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E9").Value = "  +0.14%  "

$ws.Range("D10").Formula = '="41.95"'
$ws.Range("D10").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E10").Value = "  -6.21%  "

$ws.Range("E11").Value = "  -1.81%  "

$ws.Range("D12").Formula = '="1.000"'
$ws.Range("D12").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E12").Value = "  +0.06%  "

$ws.Range("D13").Formula = '="20.62"'
$ws.Range("D13").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E13").Value = "  -4.80%  "

$ws.Range("D14").Formula = '="6.027"'
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E14").Value = "  -1.94%  "

$ws.Range("D15").Formula = '="7.122"'
$ws.Range("D15").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E15").Value = "  -2.65%  "

$ws.Range("D16").Value = "1.752.51"
$ws.Range("E16").Value = "  -1.93%  "

$ws.Range("D17").Formula = '="93.43"'
$ws.Range("D17").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E17").Value = "  +1.39%  "

$ws.Range("D18").Formula = '="0.00001060"'
$ws.Range("D18").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E18").Value = "  -0.36%  "

$ws.Range("D19").Formula = '="0.06390"'
$ws.Range("D19").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E19").Value = "  +0.63%  "

$ws.Range("D20").Formula = '="0.9996"'
$ws.Range("D20").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E20").Value = "  -0.02%  "

$ws.Range("D21").Formula = '="16.79"'
$ws.Range("D21").Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E21").Value = "  -2.58%  "

$ws.Range("D22").Formula = '="5.852"'
$ws.Range("D22").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E22").Value = "  -1.89%  "

$ws.Range("D23").Value = "27.571.85"

$ws.Range("E24").Value = "  -1.89%  "

$ws.Range("D25").Formula = '="2.107"'
$ws.Range("D25").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E25").Value = "  -0.72%  "

$ws.Range("D26").Formula = '="161.64"'
$ws.Range("D26").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E26").Value = "  +1.75%  "

$ws.Range("D27").Formula = '="20.45"'
$ws.Range("D27").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E27").Value = "  +0.35%  "

$ws.Range("D28").Value = "1.952.92"
$ws.Range("E28").Value = "  -2.04%  "

$ws.Range("D29").Formula = '="2.085"'
$ws.Range("D29").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E29").Value = "  -3.91%  "

$ws.Range("D30").Formula = '="125.28"'
$ws.Range("D30").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E30").Value = "  -1.49%  "

$ws.Range("D31").Formula = '="1.079"'
$ws.Range("D31").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E31").Value = "  -6.94%  "

$ws.Range("D32").Formula = '="3.659"'
$ws.Range("D32").Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E32").Value = "  +3.49%  "

$ws.Range("D33").Formula = '="0.09011"'
$ws.Range("D33").Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("D34").Formula = '="5.540"'
$ws.Range("D34").Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E34").Value = "  -3.15%  "

$ws.Range("D35").Formula = '="11.96"'
$ws.Range("D35").Copy() | Out-Null
$ws.Range("D35").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E35").Value = "  -5.02%  "

$ws.Range("D36").Formula = '="0.02299"'
$ws.Range("D36").Copy() | Out-Null
$ws.Range("D36").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E36").Value = "  -0.89%  "

$ws.Range("D37").Formula = '="0.06019"'
$ws.Range("D37").Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E37").Value = "  -0.49%  "

$ws.Range("D38").Formula = '="0.2088"'
$ws.Range("D38").Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E38").Value = "  -1.27%  "

$ws.Range("D39").Formula = '="0.6345"'
$ws.Range("D39").Copy() | Out-Null
$ws.Range("D39").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E39").Value = "  -1.67%  "

$ws.Range("D40").Formula = '="4.945"'
$ws.Range("D40").Copy() | Out-Null
$ws.Range("D40").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E40").Value = "  -2.90%  "

$ws.Range("D41").Formula = '="1.204"'
$ws.Range("D41").Copy() | Out-Null
$ws.Range("D41").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E41").Value = "  +1.34%  "

$ws.Range("D43").Formula = '="1.381"'
$ws.Range("D43").Copy() | Out-Null
$ws.Range("D43").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E43").Value = "  -2.77%  "

$ws.Range("D44").Formula = '="7.752"'
$ws.Range("D44").Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E44").Value = "  -1.06%  "

$ws.Range("D45").Formula = '="13.15"'
$ws.Range("D45").Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E45").Value = "  -3.58%  "

$ws.Range("D46").Formula = '="0.5884"'
$ws.Range("D46").Copy() | Out-Null
$ws.Range("D46").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E46").Value = "  -1.64%  "

$ws.Range("D47").Formula = '="3.712"'
$ws.Range("D47").Copy() | Out-Null
$ws.Range("D47").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E47").Value = "  +0.20%  "

$ws.Range("D48").Formula = '="121.93"'
$ws.Range("D48").Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E48").Value = "  -2.23%  "

$ws.Range("D49").Formula = '="1.959"'
$ws.Range("D49").Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E49").Value = "  -1.01%  "

$ws.Range("D50").Formula = '="1.147"'
$ws.Range("D50").Copy() | Out-Null
$ws.Range("D50").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E50").Value = "  -0.43%  "

$ws.Range("E51").Value = "  -1.26%  "
